$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.099562666666667
$ws.Range("H2").Value = 9.298688
$ws.Range("I2").Value = 0.2686390288432488
$ws.Range("J2").Value = 0.2686390288432488
$ws.Range("M2").Value = 19.163974
$ws.Range("N2").Value = 57.491922
$ws.Range("O2").Value = 0.6845732287637933
$ws.Range("P2").Value = 0.6845732287637933
$ws.Range("Q2").Value = 59.39993835537067
$ws.Range("R2").Value = 534.5994451983361
$ws.Range("S2").Value = 0.1839030873471926
$ws.Range("T2").Value = 0.1839030873471926
$ws.Range("G3").Value = 3.099562666666667
$ws.Range("H3").Value = 9.298688
$ws.Range("I3").Value = 0.2686390288432488
$ws.Range("J3").Value = 0.2686390288432488
$ws.Range("O3").Value = 0.02733363438148322
$ws.Range("P3").Value = 0.02733363438148323
$ws.Range("Q3").Value = 2.371720261717333
$ws.Range("R3").Value = 21.345482355456
$ws.Range("S3").Value = 0.007342880994998087
$ws.Range("T3").Value = 0.007342880994998088
$ws.Range("G4").Value = 3.099562666666667
$ws.Range("H4").Value = 9.298688
$ws.Range("I4").Value = 0.2686390288432488
$ws.Range("J4").Value = 0.2686390288432488
$ws.Range("M4").Value = 7.880893333333333
$ws.Range("N4").Value = 23.64268
$ws.Range("O4").Value = 0.281520346184098
$ws.Range("P4").Value = 0.281520346184098
$ws.Range("Q4").Value = 24.42732275598222
$ws.Range("R4").Value = 219.84590480384
$ws.Range("S4").Value = 0.07562735239851127
$ws.Range("T4").Value = 0.07562735239851127
$ws.Range("G5").Value = 3.099562666666667
$ws.Range("H5").Value = 9.298688
$ws.Range("I5").Value = 0.2686390288432488
$ws.Range("J5").Value = 0.2686390288432488
$ws.Range("M5").Value = 0.183999
$ws.Range("N5").Value = 0.551997
$ws.Range("O5").Value = 0.006572790670625477
$ws.Range("P5").Value = 0.006572790670625476
$ws.Range("Q5").Value = 0.570316431104
$ws.Range("R5").Value = 5.132847879936
$ws.Range("S5").Value = 0.001765708102546794
$ws.Range("T5").Value = 0.001765708102546794
$ws.Range("G6").Value = 6.189892666666666
$ws.Range("I6").Value = 0.5364778626674904
$ws.Range("J6").Value = 0.5364778626674905
$ws.Range("M6").Value = 19.163974
$ws.Range("N6").Value = 57.491922
$ws.Range("O6").Value = 0.6845732287637933
$ws.Range("P6").Value = 0.6845732287637933
$ws.Range("Q6").Value = 118.6229421267907
$ws.Range("R6").Value = 1067.606479141116
$ws.Range("S6").Value = 0.3672583826065828
$ws.Range("T6").Value = 0.3672583826065829
$ws.Range("G7").Value = 6.189892666666666
$ws.Range("I7").Value = 0.5364778626674904
$ws.Range("J7").Value = 0.5364778626674905
$ws.Range("O7").Value = 0.02733363438148322
$ws.Range("P7").Value = 0.02733363438148323
$ws.Range("Q7").Value = 4.736375880787333
$ws.Range("S7").Value = 0.01466388975191275
$ws.Range("T7").Value = 0.01466388975191276
$ws.Range("G8").Value = 6.189892666666666
$ws.Range("I8").Value = 0.5364778626674904
$ws.Range("J8").Value = 0.5364778626674905
$ws.Range("M8").Value = 7.880893333333333
$ws.Range("N8").Value = 23.64268
$ws.Range("O8").Value = 0.281520346184098
$ws.Range("P8").Value = 0.281520346184098
$ws.Range("Q8").Value = 48.78188385078222
$ws.Range("R8").Value = 439.03695465704
$ws.Range("S8").Value = 0.1510294336182569
$ws.Range("T8").Value = 0.1510294336182569
$ws.Range("G9").Value = 6.189892666666666
$ws.Range("I9").Value = 0.5364778626674904
$ws.Range("J9").Value = 0.5364778626674905
$ws.Range("M9").Value = 0.183999
$ws.Range("N9").Value = 0.551997
$ws.Range("O9").Value = 0.006572790670625477
$ws.Range("P9").Value = 0.006572790670625476
$ws.Range("Q9").Value = 1.138934060774
$ws.Range("R9").Value = 10.250406546966
$ws.Range("S9").Value = 0.003526156690737977
$ws.Range("T9").Value = 0.003526156690737977
$ws.Range("G10").Value = 1.888584
$ws.Range("H10").Value = 5.665752
$ws.Range("I10").Value = 0.1636835341659699
$ws.Range("J10").Value = 0.1636835341659699
$ws.Range("M10").Value = 19.163974
$ws.Range("N10").Value = 57.491922
$ws.Range("O10").Value = 0.6845732287637933
$ws.Range("P10").Value = 0.6845732287637933
$ws.Range("Q10").Value = 36.192774672816
$ws.Range("R10").Value = 325.734972055344
$ws.Range("S10").Value = 0.1120533654794667
$ws.Range("T10").Value = 0.1120533654794667
$ws.Range("G11").Value = 1.888584
$ws.Range("H11").Value = 5.665752
$ws.Range("I11").Value = 0.1636835341659699
$ws.Range("J11").Value = 0.1636835341659699
$ws.Range("O11").Value = 0.02733363438148322
$ws.Range("P11").Value = 0.02733363438148323
$ws.Range("Q11").Value = 1.445104816536
$ws.Range("R11").Value = 13.005943348824
$ws.Range("S11").Value = 0.004474065877161639
$ws.Range("T11").Value = 0.004474065877161639
$ws.Range("G12").Value = 1.888584
$ws.Range("H12").Value = 5.665752
$ws.Range("I12").Value = 0.1636835341659699
$ws.Range("J12").Value = 0.1636835341659699
$ws.Range("M12").Value = 7.880893333333333
$ws.Range("N12").Value = 23.64268
$ws.Range("O12").Value = 0.281520346184098
$ws.Range("P12").Value = 0.281520346184098
$ws.Range("Q12").Value = 14.88372905504
$ws.Range("R12").Value = 133.95356149536
$ws.Range("S12").Value = 0.04608024520304048
$ws.Range("T12").Value = 0.04608024520304049
$ws.Range("G13").Value = 1.888584
$ws.Range("H13").Value = 5.665752
$ws.Range("I13").Value = 0.1636835341659699
$ws.Range("J13").Value = 0.1636835341659699
$ws.Range("M13").Value = 0.183999
$ws.Range("N13").Value = 0.551997
$ws.Range("O13").Value = 0.006572790670625477
$ws.Range("P13").Value = 0.006572790670625476
$ws.Range("Q13").Value = 0.347497567416
$ws.Range("R13").Value = 3.127478106744
$ws.Range("S13").Value = 0.001075857606301093
$ws.Range("T13").Value = 0.001075857606301093
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3599813333333333
$ws.Range("H14").Value = 1.079944
$ws.Range("I14").Value = 0.03119957432329092
$ws.Range("J14").Value = 0.03119957432329093
$ws.Range("M14").Value = 19.163974
$ws.Range("N14").Value = 57.491922
$ws.Range("O14").Value = 0.6845732287637933
$ws.Range("P14").Value = 0.6845732287637933
$ws.Range("Q14").Value = 6.898672912485333
$ws.Range("R14").Value = 62.088056212368
$ws.Range("S14").Value = 0.02135839333055121
$ws.Range("T14").Value = 0.02135839333055121
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3599813333333333
$ws.Range("H15").Value = 1.079944
$ws.Range("I15").Value = 0.03119957432329092
$ws.Range("J15").Value = 0.03119957432329093
$ws.Range("O15").Value = 0.02733363438148322
$ws.Range("P15").Value = 0.02733363438148323
$ws.Range("Q15").Value = 0.2754501566586666
$ws.Range("R15").Value = 2.479051409928
$ws.Range("S15").Value = 0.0008527977574107459
$ws.Range("T15").Value = 0.0008527977574107461
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3599813333333333
$ws.Range("H16").Value = 1.079944
$ws.Range("I16").Value = 0.03119957432329092
$ws.Range("J16").Value = 0.03119957432329093
$ws.Range("M16").Value = 7.880893333333333
$ws.Range("N16").Value = 23.64268
$ws.Range("O16").Value = 0.281520346184098
$ws.Range("P16").Value = 0.281520346184098
$ws.Range("Q16").Value = 2.836974489991111
$ws.Range("R16").Value = 25.53277040992
$ws.Range("S16").Value = 0.008783314964289355
$ws.Range("T16").Value = 0.008783314964289357
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3599813333333333
$ws.Range("H17").Value = 1.079944
$ws.Range("I17").Value = 0.03119957432329092
$ws.Range("J17").Value = 0.03119957432329093
$ws.Range("M17").Value = 0.183999
$ws.Range("N17").Value = 0.551997
$ws.Range("O17").Value = 0.006572790670625477
$ws.Range("P17").Value = 0.006572790670625476
$ws.Range("Q17").Value = 0.06623620535199999
$ws.Range("R17").Value = 0.596125848168
$ws.Range("S17").Value = 0.0002050682710396128
$ws.Range("T17").Value = 0.0002050682710396127
